$wb = $excel.ActiveWorkbook

$urlA = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/29e3d8d038653015cf5d0610901190bb61706f5a/e2e/a.md"
$urlB = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/29e3d8d038653015cf5d0610901190bb61706f5a/e2e/b.md"

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: status column (E/F) flips from "Ready for handoff" to the
# new handed-back status for both rows.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Columns E and F widen to fit the longer status text.
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

# "Latest Target File" (I) and "Latest Handback File" (J) get populated now
# that the handback has happened.
$wsZh.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

# "Latest Handback DateTime" (K) now has a real timestamp.
$wsZh.Range("K2").Value = "2016-08-15 14:35:57"
$wsZh.Range("K3").Value = "2016-08-15 14:35:57"

# Rebuild hyperlinks in row-major order so the new "Latest Target File" link
# lands between the existing A-column links (matches relationship id order).
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $urlA, "", "", "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $urlA, "", "", "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $urlB, "", "", "b.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $urlA, "", "", "a.md")

# Columns C and J widen.
$wsZh.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZh.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$wsDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$wsDe.Range("K2").Value = "2016-08-15 14:36:11"
$wsDe.Range("K3").Value = "2016-08-15 14:36:11"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $urlA, "", "", "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $urlA, "", "", "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $urlB, "", "", "b.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $urlA, "", "", "a.md")

$wsDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDe.Columns.Item(10).ColumnWidth = 39.166666666666664

Write-Output "Report regenerated for handback."
